# Update Vis Tool + Teams Data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Games" sheet - correct the stats recorded for game 13 (row 14)
# ---------------------------------------------------------------------
$games = $wb.Worksheets.Item("Games")

$games.Range("E14").Value = 93.3
$games.Range("J14").Value = 116.8
$games.Range("M14").Value = 0.608
$games.Range("N14").Value = 12.1
$games.Range("P14").Value = 0.241
$games.Range("Q14").Value = 129.7

# ---------------------------------------------------------------------
# 2) "Games" sheet - append the result of the game that was played on
#    2024-01-xx (serial 45305) vs IND, which used to be the next
#    scheduled game on the "Next" sheet.
# ---------------------------------------------------------------------
$games.Range("A42").Value = 41
$games.Range("B42").Value = 45305
$games.Range("B42").NumberFormat = "YYYY-MM-DD"
$games.Range("C42").Value = 2
$games.Range("D42").Value = 117
$games.Range("E42").Value = 98.2
$games.Range("F42").Value = 0.718
$games.Range("G42").Value = 20.6
$games.Range("H42").Value = 14.8
$games.Range("I42").Value = 0.211
$games.Range("J42").Value = 119.2
$games.Range("K42").Value = "IND"
$games.Range("L42").Value = 109
$games.Range("M42").Value = 0.511
$games.Range("N42").Value = 11.3
$games.Range("O42").Value = 28
$games.Range("P42").Value = 0.163
$games.Range("Q42").Value = 111
$games.Range("R42").Value = 1
$games.Range("S42").Value = 1

# ---------------------------------------------------------------------
# 3) "Next" sheet - the game vs IND on 45305 has now been played, so it
#    drops off the front of the upcoming-schedule list; every remaining
#    row shifts up by one.
# ---------------------------------------------------------------------
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
